$d = $word.ActiveDocument

# Update the header date
$d.Content.Find.Execute("2026-01-06 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-07 Wednesday", 2)

# Update table cells by explicit row/column position to avoid ambiguity
# between duplicate values ("73x80=5840" appears twice).
$table = $d.Tables.Item(1)

$table.Cell(1,1).Range.Text = "40×48=1920"
$table.Cell(1,2).Range.Text = "65×19=1235"
$table.Cell(1,3).Range.Text = "60×83=4980"
$table.Cell(1,4).Range.Text = "66×41=2706"
$table.Cell(1,5).Range.Text = "36×20=720"

$table.Cell(5,1).Range.Text = "90×52=4680"
$table.Cell(5,2).Range.Text = "82×70=5740"
$table.Cell(5,3).Range.Text = "24×32=768"
$table.Cell(5,4).Range.Text = "39×31=1209"
$table.Cell(5,5).Range.Text = "28×13=364"

$table.Cell(10,1).Range.Text = "15×11=165"
$table.Cell(10,2).Range.Text = "32×67=2144"
$table.Cell(10,3).Range.Text = "88×84=7392"
$table.Cell(10,4).Range.Text = "57×57=3249"
$table.Cell(10,5).Range.Text = "66×17=1122"

$table.Cell(15,1).Range.Text = "54×34=1836"
$table.Cell(15,2).Range.Text = "71×68=4828"
$table.Cell(15,3).Range.Text = "13×90=1170"
$table.Cell(15,4).Range.Text = "91×69=6279"
$table.Cell(15,5).Range.Text = "46×24=1104"

$table.Cell(20,1).Range.Text = "89×64=5696"
$table.Cell(20,2).Range.Text = "49×88=4312"
$table.Cell(20,3).Range.Text = "89×39=3471"
$table.Cell(20,4).Range.Text = "23×41=943"
$table.Cell(20,5).Range.Text = "56×32=1792"
